$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: set the brand-new group-label strings (col A) for rows 15-23 first,
#     in final shared-string order (indices 23-28) ---
$ws.Range("A15").Value = "Stream-crude glycerol"
$ws.Range("A16").Value = "Stream-pure glycerine"
$ws.Range("A17").Value = "Stream-cellulase"
$ws.Range("A19").Value = "Pretreatment reactor system"
$ws.Range("A20").Value = "Pretreatment and saccharification"
$ws.Range("A22").Value = "Cofermenation"

# --- Step 2: update existing rows 4-14 (col B label shifts + recomputed correlations) ---
$ws.Range("B4").Value = "Cane lipid content [dry wt. %]"
$ws.Range("C4").Value = 0.3486610206024408
$ws.Range("D4").Value = 0.9803523589900942
$ws.Range("E4").Value = -0.9999971322878852
$ws.Range("F4").Value = 0.9603281124611245
$ws.Range("H4").Value = 0.3474862837234513
$ws.Range("I4").Value = -0.01894716392588655
$ws.Range("J4").Value = -0.1738375017804727
$ws.Range("K4").Value = 0.1654197721847909
$ws.Range("L4").Value = -0.02015946925437877
$ws.Range("M4").Value = -0.3808344926093796
$ws.Range("N4").Value = 0.2045579573023183
$ws.Range("P4").Value = -0.1326448818497953

$ws.Range("B5").Value = "Relative sorghum lipid content [dry wt. %]"
$ws.Range("C5").Value = 0.005349108501964339
$ws.Range("D5").Value = -0.01411938805277552
$ws.Range("E5").Value = 0.007769103478764138
$ws.Range("F5").Value = -0.002326002141040085
$ws.Range("H5").Value = -0.0004778683391147335
$ws.Range("I5").Value = -0.00162593900903756
$ws.Range("J5").Value = 0.007098927706779668
$ws.Range("K5").Value = 0.009605336928213476
$ws.Range("L5").Value = -0.01842566310502652
$ws.Range("M5").Value = 0.001891841739673669
$ws.Range("N5").Value = 0.01097233474289339
$ws.Range("P5").Value = -0.0002730643309225732

$ws.Range("B6").Value = "Lipid retention [%]"
$ws.Range("C6").Value = -0.02542045532081821
$ws.Range("D6").Value = -0.02655531120621244
$ws.Range("E6").Value = -0.02046478286659131
$ws.Range("F6").Value = 0.05895996610239863
$ws.Range("H6").Value = 0.01577479368699174
$ws.Range("I6").Value = -0.002647487913899516
$ws.Range("J6").Value = -0.008125356293979995
$ws.Range("K6").Value = -0.05681807881672314
$ws.Range("L6").Value = -0.2331732237269289
$ws.Range("M6").Value = -0.00482730652909226
$ws.Range("N6").Value = 0.1520461064018442
$ws.Range("P6").Value = 0.02308932869957314

$ws.Range("B7").Value = "Bagasse lipid extraction efficiency [%]"
$ws.Range("C7").Value = 0.1096520102100804
$ws.Range("D7").Value = 0.1652420382576815
$ws.Range("E7").Value = 0.01566537806661512
$ws.Range("F7").Value = -0.1470257862010314
$ws.Range("H7").Value = -0.0438646052425842
$ws.Range("I7").Value = -0.0005958812398352496
$ws.Range("J7").Value = 0.08039389351161819
$ws.Range("K7").Value = 0.1230824022992961
$ws.Range("L7").Value = 0.9654484080259363
$ws.Range("M7").Value = -0.007095883291835331
$ws.Range("N7").Value = -0.6601621398784855
$ws.Range("P7").Value = -0.1681691580547663

$ws.Range("B8").Value = "Capacity [ton/hr]"
$ws.Range("C8").Value = 0.2524538617141545
$ws.Range("D8").Value = -0.01699660819986433
$ws.Range("E8").Value = 0.01874287879771515
$ws.Range("F8").Value = 0.1847359375654375
$ws.Range("H8").Value = 0.9249928702477147
$ws.Range("I8").Value = 0.9576457947858317
$ws.Range("J8").Value = -0.1621446081203423
$ws.Range("K8").Value = 0.05741462418458496
$ws.Range("L8").Value = 0.006964821878592875
$ws.Range("M8").Value = -0.289063275306531
$ws.Range("N8").Value = 0.1598735103949404
$ws.Range("P8").Value = 0.4509692472867698

$ws.Range("B9").Value = "Price [USD/gal]"
$ws.Range("C9").Value = 0.5025083990763358
$ws.Range("D9").Value = -0.01448220720328829
$ws.Range("E9").Value = 0.01462552023302081
$ws.Range("F9").Value = -0.0215808599992344
$ws.Range("H9").Value = -0.03039311363172454
$ws.Range("I9").Value = -0.02565199811407992
$ws.Range("J9").Value = 0.002769786020574753
$ws.Range("K9").Value = -0.5711248544609941
$ws.Range("L9").Value = -0.0006466885698675427
$ws.Range("M9").Value = -0.003845828601833144
$ws.Range("N9").Value = 0.003873127450925097
$ws.Range("P9").Value = 0.02063873468154939

$ws.Range("B10").Value = "Price [USD/gal]"
$ws.Range("C10").Value = 0.5517150598926023
$ws.Range("D10").Value = -0.0154560975782439
$ws.Range("E10").Value = 0.01327732718709308
$ws.Range("F10").Value = -0.009758973318358932
$ws.Range("H10").Value = 0.006196915255876609
$ws.Range("I10").Value = 0.007491558251662329
$ws.Range("J10").Value = -0.004793967993768524
$ws.Range("K10").Value = 0.6141400814296032
$ws.Range("L10").Value = -0.0128048950401958
$ws.Range("M10").Value = 0.02060088197603528
$ws.Range("N10").Value = 0.007787982167519285
$ws.Range("P10").Value = 0.01352683330907333

$ws.Range("B11").Value = "Price [USD/cf]"
$ws.Range("C11").Value = -0.001537391677495667
$ws.Range("D11").Value = 0.01572554050102162
$ws.Range("E11").Value = -0.01713100388524015
$ws.Range("F11").Value = 0.01629449940377998
$ws.Range("H11").Value = 0.01231949319677973
$ws.Range("I11").Value = 0.008194044615761783
$ws.Range("J11").Value = -0.004683385397114324
$ws.Range("K11").Value = 0.01149523485980939
$ws.Range("L11").Value = -0.00182088650483546
$ws.Range("M11").Value = -0.008204314600172583
$ws.Range("N11").Value = 0.02386473973858959
$ws.Range("P11").Value = -0.004849088161963526

$ws.Range("B12").Value = "Electricity price [USD/kWh]"
$ws.Range("C12").Value = 0.08317605635104225
$ws.Range("D12").Value = -0.01032835874913435
$ws.Range("E12").Value = 0.01038483612739344
$ws.Range("F12").Value = -0.003915402684616107
$ws.Range("H12").Value = 0.01651080613243224
$ws.Range("I12").Value = 0.02000939523237581
$ws.Range("J12").Value = 0.01711309665625601
$ws.Range("K12").Value = 0.02094203805368152
$ws.Range("L12").Value = 0.005535897149435885
$ws.Range("M12").Value = -0.009627247585089902
$ws.Range("N12").Value = -0.0009434797817391911
$ws.Range("P12").Value = 0.01548816657152666

$ws.Range("B13").Value = "Operating days [day/yr]"
$ws.Range("C13").Value = 0.127034926553397
$ws.Range("D13").Value = -0.001687440835497633
$ws.Range("E13").Value = -0.0002626734825069393
$ws.Range("F13").Value = -0.0005188972047558882
$ws.Range("H13").Value = -0.01099128447165138
$ws.Range("I13").Value = 0.267955776926231
$ws.Range("J13").Value = -0.00284714906383081
$ws.Range("K13").Value = 0.009992666223706648
$ws.Range("L13").Value = -0.01599451619178064
$ws.Range("M13").Value = -0.01295039965401599
$ws.Range("N13").Value = 0.01349464681178587
$ws.Range("P13").Value = -0.006435033569401342

$ws.Range("B14").Value = "IRR [%]"
$ws.Range("C14").Value = -0.2773507782300311
$ws.Range("D14").Value = 0.01697964538318581
$ws.Range("E14").Value = -0.01727266072290643
$ws.Range("F14").Value = 0.01618241047129642
$ws.Range("H14").Value = 0.02408011209920448
$ws.Range("I14").Value = 0.02607368293094731
$ws.Range("J14").Value = -0.02340988682527531
$ws.Range("K14").Value = -0.03103775135351005
$ws.Range("L14").Value = 0.005143398733735948
$ws.Range("M14").Value = -0.02438988750359549
$ws.Range("N14").Value = 0.02063012482520499
$ws.Range("P14").Value = 0.02608065934722637

# --- Step 3: copy A4:B4 cell formatting (bold font + border, style index 1) down to rows 15-23 ---
$ws.Range("A4:B4").Copy()
$ws.Range("A15:B23").PasteSpecial(-4122)

# --- Step 4: populate new rows 15-23 (col B labels + numeric values) ---
$ws.Range("B15").Value = "Price [USD/kg]"
$ws.Range("C15").Value = 0.01468666647546666
$ws.Range("D15").Value = -0.009949870189994806
$ws.Range("E15").Value = 0.009796767559870702
$ws.Range("F15").Value = -0.00851142427645697
$ws.Range("H15").Value = -0.001734651717386069
$ws.Range("I15").Value = 0.004215266088610644
$ws.Range("J15").Value = 0.008493221925462724
$ws.Range("K15").Value = 0.0343054802202192
$ws.Range("L15").Value = 0.00385232377009295
$ws.Range("M15").Value = -0.01392426516497061
$ws.Range("N15").Value = -0.002668205866728234
$ws.Range("P15").Value = 0.003810548696421947

$ws.Range("B16").Value = "Price [USD/kg]"
$ws.Range("C16").Value = -0.04573330896533236
$ws.Range("D16").Value = 0.008334211149368446
$ws.Range("E16").Value = -0.00864499551379982
$ws.Range("F16").Value = 0.01374060842162434
$ws.Range("H16").Value = 0.02670928618837145
$ws.Range("I16").Value = 0.02145504949820198
$ws.Range("J16").Value = -0.01345020618209228
$ws.Range("K16").Value = -0.03080578846423153
$ws.Range("L16").Value = 0.001253469650138786
$ws.Range("M16").Value = -0.006936524245460969
$ws.Range("N16").Value = 0.001402618520104741
$ws.Range("P16").Value = 0.02957852470314098

$ws.Range("B17").Value = "Price [USD/kg]"
$ws.Range("C17").Value = 0.03879977665599106
$ws.Range("D17").Value = 0.03810672133226885
$ws.Range("E17").Value = -0.03859501603980064
$ws.Range("F17").Value = 0.04449396888375876
$ws.Range("H17").Value = 0.05078634663945385
$ws.Range("I17").Value = 0.04796831395073255
$ws.Range("J17").Value = -0.005035091683916938
$ws.Range("K17").Value = 0.01917543955101758
$ws.Range("L17").Value = -0.003987287199491488
$ws.Range("M17").Value = -0.01765495645019826
$ws.Range("N17").Value = 0.008499268371970733
$ws.Range("P17").Value = 0.0268805700672228

$ws.Range("B18").Value = "Cellulase loading [wt. % cellulose]"
$ws.Range("C18").Value = 0.0344067301762692
$ws.Range("D18").Value = 0.003844745145789806
$ws.Range("E18").Value = -0.00119402375976095
$ws.Range("F18").Value = 0.003441679721667188
$ws.Range("H18").Value = 0.01070070330802813
$ws.Range("I18").Value = 0.01683920947356838
$ws.Range("J18").Value = 0.01603746374698477
$ws.Range("K18").Value = -0.01139336906373476
$ws.Range("L18").Value = 0.01376556573462263
$ws.Range("M18").Value = 0.00307507826700313
$ws.Range("N18").Value = -0.02451089579643583
$ws.Range("P18").Value = -0.00721279804851192

$ws.Range("B19").Value = "Base cost [million USD]"
$ws.Range("C19").Value = 0.003274279810971192
$ws.Range("D19").Value = -0.01579782149591286
$ws.Range("E19").Value = 0.01687424400296976
$ws.Range("F19").Value = -0.01525776973031079
$ws.Range("H19").Value = -0.00348730852349234
$ws.Range("I19").Value = 0.004456985650279426
$ws.Range("J19").Value = 0.01853350556403818
$ws.Range("K19").Value = -0.008620442648817703
$ws.Range("L19").Value = 0.01314768570990743
$ws.Range("M19").Value = 0.006062885714515426
$ws.Range("N19").Value = -0.0204925349317014
$ws.Range("P19").Value = -0.004485800915432036

$ws.Range("B20").Value = "Glucose yield [%]"
$ws.Range("C20").Value = -0.004100143652005745
$ws.Range("D20").Value = -0.003273113026924521
$ws.Range("E20").Value = 0.0006928985557159422
$ws.Range("F20").Value = 0.002814760144590406
$ws.Range("H20").Value = -0.001145887149835486
$ws.Range("I20").Value = -0.001009648360385934
$ws.Range("J20").Value = 0.03247404232676749
$ws.Range("K20").Value = -0.02665623313024932
$ws.Range("L20").Value = -0.01020067672802707
$ws.Range("M20").Value = -0.007478347403133895
$ws.Range("N20").Value = -0.009736784357471374
$ws.Range("P20").Value = 0.00009312816372512654

$ws.Range("B21").Value = "Xylose yield [%]"
$ws.Range("C21").Value = 0.02471810470072418
$ws.Range("D21").Value = 0.01921914681676587
$ws.Range("E21").Value = -0.02174361562174462
$ws.Range("F21").Value = 0.02533763563750542
$ws.Range("H21").Value = 0.02899969565598782
$ws.Range("I21").Value = 0.02411104560444182
$ws.Range("J21").Value = -0.003508912371683103
$ws.Range("K21").Value = -0.01955308359812334
$ws.Range("L21").Value = -0.01101838604073544
$ws.Range("M21").Value = -0.01040513696020548
$ws.Range("N21").Value = 0.02107069649082786
$ws.Range("P21").Value = 0.0213986276079451

$ws.Range("B22").Value = "Glucose to ethanol yield [%]"
$ws.Range("C22").Value = 0.007810921656436864
$ws.Range("D22").Value = 0.006190503223620128
$ws.Range("E22").Value = -0.008760888734435548
$ws.Range("F22").Value = 0.008738446717537866
$ws.Range("H22").Value = 0.0001098464683938587
$ws.Range("I22").Value = -0.0009385170615406823
$ws.Range("J22").Value = -0.01005614250324293
$ws.Range("K22").Value = -0.01060462573618503
$ws.Range("L22").Value = -0.01427929583517183
$ws.Range("M22").Value = 0.003639984241599369
$ws.Range("N22").Value = 0.007863030938521236
$ws.Range("P22").Value = 0.001276443411057736

$ws.Range("B23").Value = "Xylose to ethanol yield [%]"
$ws.Range("C23").Value = 0.003037175641487025
$ws.Range("D23").Value = 0.003273531970941278
$ws.Range("E23").Value = -0.000710742556429702
$ws.Range("F23").Value = 0.0005396895575875822
$ws.Range("H23").Value = -0.007476486443059457
$ws.Range("I23").Value = -0.01020257877610315
$ws.Range("J23").Value = 0.01364034859009499
$ws.Range("K23").Value = 0.01179537685581507
$ws.Range("L23").Value = 0.01263250178530007
$ws.Range("M23").Value = -0.02030841873233675
$ws.Range("N23").Value = -0.01164514577780583
$ws.Range("P23").Value = -0.002042669841706793

# --- Step 5: add new merged cells for multi-row parameter group labels ---
$ws.Range("A17:A18").Merge()
$ws.Range("A20:A21").Merge()
$ws.Range("A22:A23").Merge()
